$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.621.78'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '3.426.97'
$ws.Range('E3').Value = '  +4.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '256.51'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '658.43'
$ws.Range('E6').Value = '  +5.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.50'
$ws.Range('E7').Value = '  +5.85%  '
$ws.Range('E8').Value = '  +8.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.07'
$ws.Range('E9').Value = '  +10.94%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '3.425.42'
$ws.Range('E11').Value = '  +4.26%  '
$ws.Range('E12').Value = '  +6.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.22'
$ws.Range('E13').Value = '  +6.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.38'
$ws.Range('E14').Value = '  +16.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000262'
$ws.Range('E15').Value = '  +6.06%  '
$ws.Range('D16').Value = '97.415.96'
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').Value = '4.056.74'
$ws.Range('E17').Value = '  +4.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.73'
$ws.Range('E18').Value = '  +38.27%  '
$ws.Range('D19').Value = '3.420.57'
$ws.Range('E19').Value = '  +4.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.72'
$ws.Range('E20').Value = '  +15.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.517'
$ws.Range('E21').Value = '  +62.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.98'
$ws.Range('E22').Value = '  +18.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.48'
$ws.Range('E23').Value = '  +1.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '509.74'
$ws.Range('E24').Value = '  +4.51%  '
$ws.Range('E25').Value = '  +3.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.14'
$ws.Range('E26').Value = '  +8.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '97.05'
$ws.Range('E27').Value = '  +9.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.89'
$ws.Range('E28').Value = '  +7.69%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.608.79'
$ws.Range('E29').Value = '  +5.43%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.156'
$ws.Range('E30').Value = '  +15.18%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.56'
$ws.Range('E31').Value = '  +11.68%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.201'
$ws.Range('E32').Value = '  +5.82%  '
$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.584'
$ws.Range('E34').Value = '  +23.79%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '30.15'
$ws.Range('E36').Value = '  +8.42%  '
$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.28'
$ws.Range('E37').Value = '  +18.18%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.90'
$ws.Range('E38').Value = '  +9.83%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.156'
$ws.Range('E39').Value = '  +4.62%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.43'
$ws.Range('E40').Value = '  +16.14%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '520.71'
$ws.Range('E41').Value = '  +6.40%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.71'
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.872'
$ws.Range('E43').Value = '  +13.15%  '
$ws.Range('B44').Value = 'MantraDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.70'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0419'
$ws.Range('E45').Value = '  +25.66%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.32'
$ws.Range('E46').Value = '  +8.01%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.51'
$ws.Range('E47').Value = '  +16.99%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.25'
$ws.Range('E48').Value = '  +13.77%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.59'
$ws.Range('E50').Value = '  +17.42%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.10'
$ws.Range('E51').Value = '  +7.92%  '
